$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Sheet2")
$dst = $wb.Worksheets.Item("ProgramResult")
$src.Copy($null, $src)
Write-Host "Copied"
foreach ($ws in $wb.Worksheets) {
    Write-Host $ws.Name
}
